$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.800.23"
$ws.Range("E2").Value = "  +0.31%  "
$ws.Range("D3").Value = "1.644.18"
$ws.Range("E3").Value = "  -0.19%  "
$ws.Range("E4").Value = "  +0.70%  "
$ws.Range("D5").Value = "'217.05"
$ws.Range("E5").Value = "  +0.55%  "
$ws.Range("E6").Value = "  -0.07%  "
$ws.Range("E7").Value = "  +0.54%  "
$ws.Range("D8").Value = "'0.251"
$ws.Range("E8").Value = "  -0.57%  "
$ws.Range("D9").Value = "'0.0626"
$ws.Range("E9").Value = "  -0.20%  "
$ws.Range("D10").Value = "'19.16"
$ws.Range("E10").Value = "  -0.71%  "
$ws.Range("D11").Value = "'0.0841"
$ws.Range("E11").Value = "  -0.29%  "
$ws.Range("D12").Value = "1.870.12"
$ws.Range("E12").Value = "  -0.47%  "
$ws.Range("D13").Value = "1.647.98"
$ws.Range("E13").Value = "  +0.01%  "
$ws.Range("D14").Value = "'4.17"
$ws.Range("E14").Value = "  -0.92%  "
$ws.Range("D15").Value = "'0.526"
$ws.Range("E15").Value = "  -1.30%  "
$ws.Range("D16").Value = "'64.58"
$ws.Range("E16").Value = "  -2.63%  "
$ws.Range("D17").Value = "26.815.35"
$ws.Range("E17").Value = "  +0.21%  "
$ws.Range("D18").Value = "0.0₃0737"
$ws.Range("E18").Value = "  -2.09%  "
$ws.Range("D19").Value = "'214.37"
$ws.Range("E19").Value = "  -2.37%  "
$ws.Range("E20").Value = "  +0.59%  "
$ws.Range("E21").Value = "  -0.50%  "
$ws.Range("D22").Value = "'2.42"
$ws.Range("E22").Value = "  +13.92%  "
$ws.Range("D23").Value = "'6.29"
$ws.Range("E23").Value = "  -0.64%  "
$ws.Range("D24").Value = "'9.36"
$ws.Range("E24").Value = "  -1.97%  "
$ws.Range("D25").Value = "'145.04"
$ws.Range("E25").Value = "  -1.18%  "
$ws.Range("E26").Value = "  +0.50%  "
$ws.Range("D27").Value = "'0.118"
$ws.Range("E27").Value = "  -1.91%  "
$ws.Range("D28").Value = "'7.09"
$ws.Range("E28").Value = "  +0.07%  "
$ws.Range("D29").Value = "'15.68"
$ws.Range("E29").Value = "  -1.04%  "
$ws.Range("D30").Value = "'0.0515"
$ws.Range("E30").Value = "  -0.69%  "
$ws.Range("E31").Value = "  +0.48%  "
$ws.Range("D32").Value = "'3.32"
$ws.Range("E32").Value = "  -2.87%  "
$ws.Range("E33").Value = "  -1.32%  "
$ws.Range("D34").Value = "1.286.06"
$ws.Range("E34").Value = "  +0.10%  "
$ws.Range("E35").Value = "  -1.15%  "
$ws.Range("E36").Value = "  +1.21%  "
$ws.Range("E37").Value = "  -4.92%  "
$ws.Range("D38").Value = "'0.539"
$ws.Range("E38").Value = "  +3.00%  "
$ws.Range("D39").Value = "'0.825"
$ws.Range("E39").Value = "  -0.34%  "
$ws.Range("E40").Value = "  +0.57%  "
$ws.Range("D41").Value = "'0.810"
$ws.Range("E41").Value = "  -0.14%  "
$ws.Range("D43").Value = "'5.35"
$ws.Range("E43").Value = "  -1.40%  "
$ws.Range("D44").Value = "1.796.28"
$ws.Range("E44").Value = "  +0.37%  "
$ws.Range("D45").Value = "'91.42"
$ws.Range("E45").Value = "  -2.30%  "
$ws.Range("D46").Value = "'60.24"
$ws.Range("E46").Value = "  +1.55%  "
$ws.Range("D47").Value = "'1.60"
$ws.Range("E47").Value = "  -0.20%  "
$ws.Range("B48").Value = "Cronos"
$ws.Range("C48").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D48").Value = "'0.0521"
$ws.Range("E48").Value = "  +0.79%  "
$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").Value = "'7.70"
$ws.Range("E49").Value = "  -0.69%  "
$ws.Range("B50").Value = "Algorand"
$ws.Range("C50").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D50").Value = "'0.0978"
$ws.Range("E50").Value = "  +0.36%  "
$ws.Range("B51").Value = "Mantle"
$ws.Range("C51").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D51").Value = "'0.408"
$ws.Range("E51").Value = "  +0.09%  "
